$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27276
$ws.Range("E2").Value = 528639277835
$ws.Range("F2").Value = 13874063711
$ws.Range("G2").Value = 1.74562
$ws.Range("D3").Value = 1852.87
$ws.Range("E3").Value = 222850581871
$ws.Range("F3").Value = 7367478450
$ws.Range("G3").Value = 2.11431
$ws.Range("F4").Value = 17897337097
$ws.Range("G4").Value = 0.07394000000000001
$ws.Range("D5").Value = 313.93
$ws.Range("E5").Value = 49568226989
$ws.Range("F5").Value = 476025639
$ws.Range("G5").Value = 1.26289
$ws.Range("F6").Value = 3961433447
$ws.Range("G6").Value = 0.18388
$ws.Range("D7").Value = 0.460605
$ws.Range("F7").Value = 901457129
$ws.Range("G7").Value = 0.92757
$ws.Range("D8").Value = 0.370853
$ws.Range("E8").Value = 12994617048
$ws.Range("F8").Value = 157547707
$ws.Range("G8").Value = 0.74859
$ws.Range("D9").Value = 1851.46
$ws.Range("F9").Value = 7007026
$ws.Range("G9").Value = 2.07822
$ws.Range("D10").Value = 0.072948
$ws.Range("E10").Value = 10174200775
$ws.Range("F10").Value = 200840393
$ws.Range("G10").Value = 0.32146
$ws.Range("D11").Value = 0.888432
$ws.Range("E11").Value = 8244403754
$ws.Range("F11").Value = 214674809
$ws.Range("G11").Value = 2.40271
$ws.Range("F12").Value = 280395099
$ws.Range("G12").Value = 2.53619
$ws.Range("D13").Value = 0.078044
$ws.Range("F13").Value = 384180992
$ws.Range("G13").Value = -1.13053
$ws.Range("D14").Value = 91.36
$ws.Range("E14").Value = 6667210013
$ws.Range("F14").Value = 380655855
$ws.Range("G14").Value = 0.77098
$ws.Range("F15").Value = 96839308
$ws.Range("G15").Value = 1.66126
$ws.Range("E16").Value = 5435876053
$ws.Range("F16").Value = 1785248241
$ws.Range("G16").Value = 0.15938
$ws.Range("F17").Value = 151672461
$ws.Range("G17").Value = 1.03469
$ws.Range("D18").Value = 14.73
$ws.Range("E18").Value = 4925729567
$ws.Range("F18").Value = 126639848
$ws.Range("G18").Value = 0.89701
$ws.Range("E19").Value = 4649942367
$ws.Range("F19").Value = 98405498
$ws.Range("G19").Value = 0.12258
$ws.Range("D20").Value = 27321
$ws.Range("F20").Value = 68297836
$ws.Range("G20").Value = 1.70357
$ws.Range("D21").Value = 5.11
$ws.Range("F21").Value = 36285067
$ws.Range("G21").Value = 0.87341
$ws.Range("E22").Value = 3373022006
$ws.Range("F22").Value = 106052224
$ws.Range("G22").Value = 0.269
$ws.Range("E23").Value = 3284950488
$ws.Range("F23").Value = 249464
$ws.Range("G23").Value = 0.50552
$ws.Range("D24").Value = 10.54
$ws.Range("E24").Value = 3084062596
$ws.Range("F24").Value = 71047397
$ws.Range("G24").Value = 0.76315
$ws.Range("F25").Value = 10497130
$ws.Range("G25").Value = 5.24755
$ws.Range("D26").Value = 46.14
$ws.Range("F26").Value = 5840252
$ws.Range("G26").Value = 0.7467200000000001
$ws.Range("D27").Value = 151.9
$ws.Range("E27").Value = 2756579606
$ws.Range("F27").Value = 71039057
$ws.Range("G27").Value = 0.01706
$ws.Range("D28").Value = 18.39
$ws.Range("E28").Value = 2595079537
$ws.Range("F28").Value = 64203204
$ws.Range("G28").Value = 1.25849
$ws.Range("D29").Value = 0.08813500000000001
$ws.Range("F29").Value = 34682590
$ws.Range("G29").Value = -0.02373
$ws.Range("D30").Value = 115.59
$ws.Range("E30").Value = 2242335298
$ws.Range("F30").Value = 53530295
$ws.Range("G30").Value = 1.27397
$ws.Range("D31").Value = 5.07
$ws.Range("E31").Value = 2209821934
$ws.Range("F31").Value = 20013311
$ws.Range("G31").Value = -0.79451
$ws.Range("F32").Value = 225910435
$ws.Range("G32").Value = 0.14667
$ws.Range("E33").Value = 1918538847
$ws.Range("F33").Value = 85299154
$ws.Range("G33").Value = 1.98419
$ws.Range("E34").Value = 1812480864
$ws.Range("F34").Value = 38836466
$ws.Range("G34").Value = 1.34509
$ws.Range("D35").Value = 0.052826
$ws.Range("E35").Value = 1662170418
$ws.Range("F35").Value = 14904065
$ws.Range("G35").Value = 0.9819
$ws.Range("E36").Value = 1658363484
$ws.Range("F36").Value = 46827604
$ws.Range("G36").Value = 3.23461
$ws.Range("D37").Value = 0.062121
$ws.Range("E37").Value = 1569624329
$ws.Range("F37").Value = 5135880
$ws.Range("G37").Value = 0.05143
$ws.Range("E38").Value = 1491076358
$ws.Range("F38").Value = 151936138
$ws.Range("G38").Value = 4.18773
$ws.Range("D39").Value = 102.46
$ws.Range("F39").Value = 12444136
$ws.Range("G39").Value = 0.70685
$ws.Range("E40").Value = 1489716781
$ws.Range("F40").Value = 51352474
$ws.Range("G40").Value = 1.61429
$ws.Range("D41").Value = 0.01952377
$ws.Range("F41").Value = 33348086
$ws.Range("G41").Value = 1.08415
$ws.Range("E42").Value = 1297243483
$ws.Range("F42").Value = 40646797
$ws.Range("G42").Value = 2.69316
$ws.Range("D43").Value = 0.164273
$ws.Range("E43").Value = 1189504884
$ws.Range("F43").Value = 30106736
$ws.Range("G43").Value = 1.06676
$ws.Range("D44").Value = 0.122792
$ws.Range("F44").Value = 28799680
$ws.Range("G44").Value = 2.44417
$ws.Range("D45").Value = 0.367123
$ws.Range("E45").Value = 1023332918
$ws.Range("F45").Value = 84048080
$ws.Range("G45").Value = 1.62506
$ws.Range("D46").Value = 1
$ws.Range("F46").Value = 41330700
$ws.Range("G46").Value = 0.06762
$ws.Range("D47").Value = 2.75
$ws.Range("F47").Value = 154827769
$ws.Range("G47").Value = 14.03482
$ws.Range("D48").Value = 1.001
$ws.Range("E48").Value = 1004791651
$ws.Range("F48").Value = 4200723
$ws.Range("G48").Value = 0.29571
$ws.Range("E49").Value = 1000291325
$ws.Range("F49").Value = 5489960
$ws.Range("G49").Value = 0.09003
$ws.Range("D50").Value = 0.513279
$ws.Range("F50").Value = 65460804
$ws.Range("G50").Value = 0.39961
$ws.Range("D51").Value = 0.859612
$ws.Range("E51").Value = 950591945
$ws.Range("F51").Value = 87428329
$ws.Range("G51").Value = 1.05194

Write-Host "Updated 162 cells across rows 2-51"
